$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the account rows: Binance moves to row 2, Bakery to row 3, Pancake to row 4
$ws.Range("A2").Value = "Binance"
$ws.Range("B2").Value = 100
$ws.Range("A3").Value = "Bakery"
$ws.Range("B3").Value = 20
$ws.Range("A4").Value = "Pancake"
$ws.Range("B4").Value = 50

# Touch the page setup so it gets serialized with explicit paper size / orientation
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Update the active selection to D7 (as if the user clicked there after adding new account)
$ws.Range("D7").Select()
